$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "89.112.52"
$ws.Range("E2").Value = "  +9.70%  "

$ws.Range("D3").Value = "3.384.82"
$ws.Range("E3").Value = "  +7.64%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "223.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.16%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "654.17"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.42%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.421"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +48.02%  "

$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.658"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +12.82%  "

$ws.Range("B9").Value = "USDC"
$ws.Range("C9").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.00%  "

$ws.Range("D10").Value = "3.380.35"
$ws.Range("E10").Value = "  +7.60%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.638"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +9.24%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000291"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +15.23%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.36"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +19.05%  "

$ws.Range("E14").Value = "  +2.85%  "

$ws.Range("B15").Value = "Toncoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.65"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.02%  "

$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "4.005.65"
$ws.Range("E16").Value = "  +7.67%  "

$ws.Range("D17").Value = "88.911.76"
$ws.Range("E17").Value = "  +9.83%  "

$ws.Range("D18").Value = "3.376.91"
$ws.Range("E18").Value = "  +7.63%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "15.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.51%  "

$ws.Range("E20").Value = "  -0.54%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +9.13%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "464.63"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.55%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.75"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +12.57%  "

$ws.Range("E24").Value = "  +3.45%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +19.45%  "

$ws.Range("D27").Value = "3.570.91"
$ws.Range("E27").Value = "  +8.05%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "81.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.57%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000145"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +19.32%  "

$ws.Range("E30").Value = "  -0.02%  "

$ws.Range("E31").Value = "  +42.07%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.49"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.93%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "594.19"
$ws.Range("D33").Style = "Normal"

$ws.Range("E34").Value = "  +0.38%  "

$ws.Range("E35").Value = "  +5.14%  "

$ws.Range("E36").Value = "  +7.23%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.34"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +22.65%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.146"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.38%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "23.99"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.69%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.439"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.05%  "

$ws.Range("E41").Value = "  +7.19%  "

$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.22"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.66%  "

$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.85"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.44%  "

$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.997"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.18%  "

$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").Style = "Normal"

$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "157.62"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.60%  "

$ws.Range("E47").Value = "  +10.27%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "189.76"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.92%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "46.88"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.60%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.57"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.29%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.676"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +8.50%  "
